# powerpoints use jsfiddle instead of others
#
# 1) The date placeholder cached on the layout behind slide 11 re-caches
#    to the new save date (8/30/2024 -> 9/4/2024).
# 2) On slide 11, the single "https://glitch.com/..." run is replaced by
#    a hyperlinked "https://jsfiddle.net/mxtdvL10/" run followed by a
#    plain run containing just a trailing space.

$p = $ppt.ActivePresentation

# "Example" slide holding the glitch.com / jsfiddle link.
$s11 = $p.Slides.Item(11)

# --- 1. Re-cache the date placeholder text on slide 11's layout ---
$layout = $s11.CustomLayout
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $lshp = $layout.Shapes.Item($i)
    if ($lshp.Name -like "Date Placeholder*") {
        $lshp.TextFrame.TextRange.Text = "9/4/2024"
    }
}

# --- 2. Swap the glitch.com link for a jsfiddle.net hyperlink ---
$newUrl = "https://jsfiddle.net/mxtdvL10/"

$contentShp = $null
for ($i = 1; $i -le $s11.Shapes.Count; $i++) {
    $candidate = $s11.Shapes.Item($i)
    if ($candidate.Name -like "Content Placeholder*") {
        $contentShp = $candidate
    }
}

$tr = $contentShp.TextFrame.TextRange

# Clear the existing run entirely (this also drops the stale endParaRPr
# echo of the old formatting) and retype the replacement text, which
# keeps the paragraph's current character formatting (bold, size, fill,
# highlight) for the freshly typed runs.
$tr.Delete() | Out-Null
$tr.InsertAfter($newUrl + " ") | Out-Null

# First run (the URL) becomes a hyperlink; the trailing space run stays
# plain, matching the two-run split in the authored slide.
$urlRange = $tr.Characters(1, $newUrl.Length)
$urlRange.ActionSettings(1).Hyperlink.Address = $newUrl
